$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing header cell (H1, style index 1: bold,
# bordered, centered) onto the two new header cells so I1/J1 match the
# other header cells exactly, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells for row 2.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
